$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 496
$ws.Range("I38").Value = 145.6
$ws.Range("K38").Value = 436.8
$ws.Range("M38").Value = -64.79999999999995
$ws.Range("H43").Value = 2274
$ws.Range("J43").Value = 2036
$ws.Range("L43").Value = 2036
$ws.Range("N43").Value = -2174
$ws.Range("H70").Value = 20086.875
$ws.Range("I70").Value = 1846
$ws.Range("J70").Value = 26167.166
$ws.Range("K70").Value = 5538
$ws.Range("L70").Value = 78501.49800000001
$ws.Range("M70").Value = -5268
$ws.Range("N70").Value = -79041.49800000001
$ws.Range("H73").Value = 20086.875
$ws.Range("I73").Value = 1846
$ws.Range("J73").Value = 26167.166
$ws.Range("K73").Value = 5538
$ws.Range("L73").Value = 78501.49800000001
$ws.Range("M73").Value = -4602
$ws.Range("N73").Value = -80373.49800000001
$ws.Range("H106").Value = 2566.5
$ws.Range("I106").Value = 3449.5
$ws.Range("J106").Value = 2125
$ws.Range("K106").Value = 3449.5
$ws.Range("L106").Value = 2125
$ws.Range("M106").Value = -2818.5
$ws.Range("N106").Value = -3387
$ws.Range("H111").Value = 4236.75
$ws.Range("I111").Value = 3599.5
$ws.Range("K111").Value = 10798.5
$ws.Range("M111").Value = -7731.5
$ws.Range("H112").Value = 1711.8
$ws.Range("J112").Value = 1792.3125
$ws.Range("L112").Value = 5376.9375
$ws.Range("N112").Value = -7592.9375
$ws.Range("H115").Value = 525.9375
$ws.Range("I115").Value = 525.9375
$ws.Range("K115").Value = 1577.8125
$ws.Range("M115").Value = -10.8125
$ws.Range("H137").Value = 1470610.2
$ws.Range("I137").Value = 1158700.6
$ws.Range("J137").Value = 1648844.4
$ws.Range("K137").Value = 3476101.8
$ws.Range("L137").Value = 4946533.199999999
$ws.Range("M137").Value = -3473551.8
$ws.Range("N137").Value = -4951633.199999999
$ws.Range("H138").Value = 3360.5557
$ws.Range("I138").Value = 1801.174
$ws.Range("J138").Value = 4517.516
$ws.Range("K138").Value = 5403.522
$ws.Range("L138").Value = 13552.548
$ws.Range("M138").Value = -263.5219999999999
$ws.Range("N138").Value = -23832.548
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1342163.5
$ws.Range("I61").Value = 1458634.4
$ws.Range("J61").Value = 2748.5
$ws.Range("K61").Value = 1458634.4
$ws.Range("L61").Value = 2748.5
$ws.Range("M61").Value = -1458422.4
$ws.Range("N61").Value = -3172.5
$ws.Range("H62").Value = 56281.125
$ws.Range("J62").Value = 56281.125
$ws.Range("L62").Value = 56281.125
$ws.Range("N62").Value = -57529.125
$ws.Range("H65").Value = 56281.125
$ws.Range("J65").Value = 56281.125
$ws.Range("L65").Value = 168843.375
$ws.Range("N65").Value = -175083.375
$ws.Range("H74").Value = 3679808.8
$ws.Range("I74").Value = 4632080
$ws.Range("J74").Value = 6762.857
$ws.Range("K74").Value = 4632080
$ws.Range("L74").Value = 6762.857
$ws.Range("M74").Value = -4631206
$ws.Range("N74").Value = -8510.857
$ws.Range("H77").Value = 3679808.8
$ws.Range("I77").Value = 4632080
$ws.Range("J77").Value = 6762.857
$ws.Range("K77").Value = 23160400
$ws.Range("L77").Value = 33814.285
$ws.Range("M77").Value = -23156032
$ws.Range("N77").Value = -42550.285
$ws.Range("H98").Value = 80349.664
$ws.Range("J98").Value = 80349.664
$ws.Range("L98").Value = 80349.664
$ws.Range("N98").Value = -86339.664
$ws.Range("H102").Value = 3065.8462
$ws.Range("I102").Value = 2360.5908
$ws.Range("J102").Value = 6944.75
$ws.Range("K102").Value = 2360.5908
$ws.Range("L102").Value = 6944.75
$ws.Range("M102").Value = -738.5907999999999
$ws.Range("N102").Value = -10188.75
$ws.Range("H132").Value = 1015066.8
$ws.Range("I132").Value = 1369867
$ws.Range("K132").Value = 4109601
$ws.Range("M132").Value = -4107071
$ws.Range("H136").Value = 1342163.5
$ws.Range("I136").Value = 1458634.4
$ws.Range("J136").Value = 2748.5
$ws.Range("K136").Value = 4375903.199999999
$ws.Range("L136").Value = 8245.5
$ws.Range("M136").Value = -4373353.199999999
$ws.Range("N136").Value = -13345.5
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1948.5714
$ws.Range("I86").Value = 2080.7334
$ws.Range("J86").Value = 1618.1666
$ws.Range("K86").Value = 2080.7334
$ws.Range("L86").Value = 1618.1666
$ws.Range("M86").Value = -957.7334000000001
$ws.Range("N86").Value = -3864.1666
$ws.Range("H89").Value = 1948.5714
$ws.Range("I89").Value = 2080.7334
$ws.Range("J89").Value = 1618.1666
$ws.Range("K89").Value = 10403.667
$ws.Range("L89").Value = 8090.833000000001
$ws.Range("M89").Value = -4787.667000000001
$ws.Range("N89").Value = -19322.833
$ws.Range("H105").Value = 1858.05
$ws.Range("I105").Value = 1924.4286
$ws.Range("J105").Value = 1703.1666
$ws.Range("K105").Value = 1924.4286
$ws.Range("L105").Value = 1703.1666
$ws.Range("M105").Value = -177.4286
$ws.Range("N105").Value = -5197.1666
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1032.625
$ws.Range("I22").Value = 973.7857
$ws.Range("K22").Value = 973.7857
$ws.Range("M22").Value = -623.7857
$ws.Range("H41").Value = 28529.375
$ws.Range("J41").Value = 37440.5
$ws.Range("L41").Value = 37440.5
$ws.Range("N41").Value = -38296.5
$ws.Range("H50").Value = 70333
$ws.Range("J50").Value = 70333
$ws.Range("L50").Value = 70333
$ws.Range("N50").Value = -71583
$ws.Range("H51").Value = 45124.625
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 45124.625
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 45124.625
$ws.Range("M51").ClearContents()
$ws.Range("N51").Value = -46596.625
$ws.Range("H59").Value = 70682.60000000001
$ws.Range("J59").Value = 70682.60000000001
$ws.Range("L59").Value = 70682.60000000001
$ws.Range("N59").Value = -72972.60000000001
$ws.Range("H60").Value = 28142.428
$ws.Range("J60").Value = 26249.5
$ws.Range("L60").Value = 26249.5
$ws.Range("N60").Value = -27271.5
$ws.Range("H61").Value = 45124.625
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 45124.625
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 45124.625
$ws.Range("M61").ClearContents()
$ws.Range("N61").Value = -45820.625
$ws.Range("H63").Value = 200080220
$ws.Range("I63").Value = 1000000000
$ws.Range("K63").Value = 1000000000
$ws.Range("M63").Value = -999999314
$ws.Range("H66").Value = 200080220
$ws.Range("I66").Value = 1000000000
$ws.Range("K66").Value = 3000000000
$ws.Range("M66").Value = -2999996568
$ws.Range("H74").Value = 47150.066
$ws.Range("J74").Value = 47333.285
$ws.Range("L74").Value = 47333.285
$ws.Range("N74").Value = -49081.285
$ws.Range("H77").Value = 47150.066
$ws.Range("J77").Value = 47333.285
$ws.Range("L77").Value = 141999.855
$ws.Range("N77").Value = -150735.855
$ws.Range("H122").Value = 2640.44
$ws.Range("I122").Value = 1199.75
$ws.Range("K122").Value = 3599.25
$ws.Range("M122").Value = -1149.25
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 2733549
$ws.Range("I7").Value = 1818444.1
$ws.Range("J7").Value = 5250087.5
$ws.Range("K7").Value = 5455332.300000001
$ws.Range("L7").Value = 15750262.5
$ws.Range("M7").Value = -5455220.300000001
$ws.Range("N7").Value = -15750486.5
$ws.Range("H116").Value = 2469.5
$ws.Range("I116").Value = 2803.4
$ws.Range("K116").Value = 8410.200000000001
$ws.Range("M116").Value = -4968.200000000001
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 114455.44
$ws.Range("J80").Value = 6704.316
$ws.Range("L80").Value = 6704.316
$ws.Range("N80").Value = -8700.315999999999
$ws.Range("H83").Value = 114455.44
$ws.Range("J83").Value = 6704.316
$ws.Range("L83").Value = 33521.58
$ws.Range("N83").Value = -43505.58
$ws.Range("H132").Value = 1006351.56
$ws.Range("I132").Value = 1340490.2
$ws.Range("K132").Value = 4021470.6
$ws.Range("M132").Value = -4018940.6
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H103").Value = 102498.5
$ws.Range("J103").Value = 102498.5
$ws.Range("L103").Value = 102498.5
$ws.Range("N103").Value = -104842.5
$ws.Range("H136").Value = 4969.6943
$ws.Range("I136").Value = 3842.7307
$ws.Range("K136").Value = 11528.1921
$ws.Range("M136").Value = -8978.1921
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 9298242
$ws.Range("I136").Value = 10891285
$ws.Range("J136").Value = 5491.6665
$ws.Range("K136").Value = 32673855
$ws.Range("L136").Value = 16474.9995
$ws.Range("M136").Value = -32671305
$ws.Range("N136").Value = -21574.9995
